# Generate Report for Handback
# Marks the two content rows of the zh-cn and de-de sheets as "handed back",
# fills in the Latest Target File (E) / Latest Handback File (F) columns
# (mirroring the source file / handoff file already in A / C), and records
# the handback timestamp in the Latest Handback DateTime column (G).

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
# RGB(100,149,237) == 0x6495ED, encoded the way VBA's RGB()/Font.Color wants
# it (so it round-trips to the existing "HyperLink" look used by A/C).
$linkColor = 15570276

# ================= zh-cn sheet =================
$ws = $wb.Worksheets.Item("zh-cn")

# --- Row 2 (80ceb9b9-...md) ---
$ws.Range("B2").Value = $statusText

$cell = $ws.Range("E2")
$ws.Hyperlinks.Add($cell, "https://github.com/OpenLocalizationTest/oltest/blob/216ff91f21d6d948791fd4958deb35d6f6b35f64/e2e/80ceb9b9-d06f-4040-b9f6-55dbbbfa287d.md", $null, $null, "80ceb9b9-d06f-4040-b9f6-55dbbbfa287d.md") | Out-Null
$cell.Font.Color = $linkColor
$cell.Font.Underline = $true
$cell.Font.Name = "Calibri"
$cell.Font.Size = 11

$cell = $ws.Range("F2")
$ws.Hyperlinks.Add($cell, "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/61a3910a74b117a483bd1ec607c5e50200cbd199/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/80ceb9b9-d06f-4040-b9f6-55dbbbfa287d.24877f15d3eefb0c26507cc2860096079cab7782.zh-cn.xlf", $null, $null, "80ceb9b9-d06f-4040-b9f6-55dbbbfa287d.24877f15d3eefb0c26507cc2860096079cab7782.zh-cn.xlf") | Out-Null
$cell.Font.Color = $linkColor
$cell.Font.Underline = $true
$cell.Font.Name = "Calibri"
$cell.Font.Size = 11

$ws.Range("G2").Value = "2016-03-04 06:42:13"

# --- Row 3 (f397695c-...md) ---
$ws.Range("B3").Value = $statusText

$cell = $ws.Range("E3")
$ws.Hyperlinks.Add($cell, "https://github.com/OpenLocalizationTest/oltest/blob/216ff91f21d6d948791fd4958deb35d6f6b35f64/e2e/f397695c-0ec3-4d6a-a189-1b31f3b14007.md", $null, $null, "f397695c-0ec3-4d6a-a189-1b31f3b14007.md") | Out-Null
$cell.Font.Color = $linkColor
$cell.Font.Underline = $true
$cell.Font.Name = "Calibri"
$cell.Font.Size = 11

$cell = $ws.Range("F3")
$ws.Hyperlinks.Add($cell, "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/61a3910a74b117a483bd1ec607c5e50200cbd199/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f397695c-0ec3-4d6a-a189-1b31f3b14007.d03e07658630811d43ed625d04768daa097b053c.zh-cn.xlf", $null, $null, "f397695c-0ec3-4d6a-a189-1b31f3b14007.d03e07658630811d43ed625d04768daa097b053c.zh-cn.xlf") | Out-Null
$cell.Font.Color = $linkColor
$cell.Font.Underline = $true
$cell.Font.Name = "Calibri"
$cell.Font.Size = 11

$ws.Range("G3").Value = "2016-03-04 06:42:13"

# ================= de-de sheet =================
$ws = $wb.Worksheets.Item("de-de")

# --- Row 2 (80ceb9b9-...md) ---
$ws.Range("B2").Value = $statusText

$cell = $ws.Range("E2")
$ws.Hyperlinks.Add($cell, "https://github.com/OpenLocalizationTest/oltest/blob/216ff91f21d6d948791fd4958deb35d6f6b35f64/e2e/80ceb9b9-d06f-4040-b9f6-55dbbbfa287d.md", $null, $null, "80ceb9b9-d06f-4040-b9f6-55dbbbfa287d.md") | Out-Null
$cell.Font.Color = $linkColor
$cell.Font.Underline = $true
$cell.Font.Name = "Calibri"
$cell.Font.Size = 11

$cell = $ws.Range("F2")
$ws.Hyperlinks.Add($cell, "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f7b0b0d69671de927707c608c01b687518645f04/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/80ceb9b9-d06f-4040-b9f6-55dbbbfa287d.24877f15d3eefb0c26507cc2860096079cab7782.de-de.xlf", $null, $null, "80ceb9b9-d06f-4040-b9f6-55dbbbfa287d.24877f15d3eefb0c26507cc2860096079cab7782.de-de.xlf") | Out-Null
$cell.Font.Color = $linkColor
$cell.Font.Underline = $true
$cell.Font.Name = "Calibri"
$cell.Font.Size = 11

$ws.Range("G2").Value = "2016-03-04 06:42:42"

# --- Row 3 (f397695c-...md) ---
$ws.Range("B3").Value = $statusText

$cell = $ws.Range("E3")
$ws.Hyperlinks.Add($cell, "https://github.com/OpenLocalizationTest/oltest/blob/216ff91f21d6d948791fd4958deb35d6f6b35f64/e2e/f397695c-0ec3-4d6a-a189-1b31f3b14007.md", $null, $null, "f397695c-0ec3-4d6a-a189-1b31f3b14007.md") | Out-Null
$cell.Font.Color = $linkColor
$cell.Font.Underline = $true
$cell.Font.Name = "Calibri"
$cell.Font.Size = 11

$cell = $ws.Range("F3")
$ws.Hyperlinks.Add($cell, "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f7b0b0d69671de927707c608c01b687518645f04/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f397695c-0ec3-4d6a-a189-1b31f3b14007.d03e07658630811d43ed625d04768daa097b053c.de-de.xlf", $null, $null, "f397695c-0ec3-4d6a-a189-1b31f3b14007.d03e07658630811d43ed625d04768daa097b053c.de-de.xlf") | Out-Null
$cell.Font.Color = $linkColor
$cell.Font.Underline = $true
$cell.Font.Name = "Calibri"
$cell.Font.Size = 11

$ws.Range("G3").Value = "2016-03-04 06:42:42"

Write-Output "Handback report generated for zh-cn and de-de."
